# Update the workbook window width (cosmetic UI state captured when the
# author resized their Excel window before saving).
$excel.Windows.Item(1).Width = 1440

# Navigate to the "Plan" worksheet, which holds the study-plan tracker.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")
$ws.Activate()

# Record actual completed hours ("实际完成") for several tasks in Stage 1
# and Stage 2 of the plan. Columns: G = standard hours, I = actual hours.
$ws.Cells.Item(17, 9).Value = 4   # I17 - task 1.1,  fully completed (G17=4)
$ws.Cells.Item(18, 9).Value = 2   # I18 - task 1.11, fully completed (G18=2)
$ws.Cells.Item(19, 9).Value = 1   # I19 - task 1.12, fully completed (G19=1)

$ws.Cells.Item(29, 9).Value = 2   # I29 - task 2.01, fully completed (G29=2)
$ws.Cells.Item(30, 9).Value = 2   # I30 - task 2.02, fully completed (G30=2)
$ws.Cells.Item(31, 9).Value = 3   # I31 - task 2.03, fully completed (G31=3)
$ws.Cells.Item(33, 9).Value = 3   # I33 - task 2.05, fully completed (G33=3)
$ws.Cells.Item(35, 9).Value = 3   # I35 - task 2.07, fully completed (G35=3)
$ws.Cells.Item(37, 9).Value = 3   # I37 - task 2.09, fully completed (G37=3)

# Move the active selection to J38 to match where the author left off.
$ws.Range("J38").Select()
